$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 544.3077
$ws.Range("I4").Value = 330.66666
$ws.Range("K4").Value = 330.66666
$ws.Range("M4").Value = -216.66666

$ws.Range("H15").Value = 1371.1666
$ws.Range("I15").Value = 1371.1666
$ws.Range("K15").Value = 4113.4998
$ws.Range("M15").Value = -3944.4998

$ws.Range("H19").Value = 1130.0834
$ws.Range("I19").Value = 1061.125
$ws.Range("K19").Value = 1061.125
$ws.Range("M19").Value = -886.125

$ws.Range("H42").Value = 3373.818
$ws.Range("I42").Value = 701.2
$ws.Range("K42").Value = 2103.6
$ws.Range("M42").Value = -1873.6

$ws.Range("H106").Value = 3431.8462
$ws.Range("I106").Value = 3431.8462
$ws.Range("K106").Value = 3431.8462
$ws.Range("M106").Value = -2800.8462

$ws.Range("H111").Value = 2089.5715
$ws.Range("I111").Value = 1417.5714
$ws.Range("K111").Value = 4252.7142
$ws.Range("M111").Value = -1185.7142

$ws.Range("H132").Value = 2923.8594
$ws.Range("I132").Value = 2388.724
$ws.Range("J132").Value = 8096.8335
$ws.Range("K132").Value = 7166.172
$ws.Range("L132").Value = 24290.5005
$ws.Range("M132").Value = -4636.172
$ws.Range("N132").Value = -29350.5005

$ws.Range("H137").Value = 4779.1055
$ws.Range("I137").Value = 1291.375
$ws.Range("J137").Value = 7315.636
$ws.Range("K137").Value = 3874.125
$ws.Range("L137").Value = 21946.908
$ws.Range("M137").Value = -1324.125
$ws.Range("N137").Value = -27046.908

$ws.Range("H138").Value = 3568.0212
$ws.Range("I138").Value = 3125.1
$ws.Range("J138").Value = 3687.7297
$ws.Range("K138").Value = 9375.299999999999
$ws.Range("L138").Value = 11063.1891
$ws.Range("M138").Value = -4235.299999999999
$ws.Range("N138").Value = -21343.1891

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10667.872
$ws.Range("I32").Value = 6760.9165
$ws.Range("K32").Value = 6760.9165
$ws.Range("M32").Value = -6473.9165

$ws.Range("H36").Value = 2262.5
$ws.Range("I36").Value = 2262.5
$ws.Range("K36").Value = 2262.5
$ws.Range("M36").Value = -1916.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5045.909
$ws.Range("I105").Value = 4919.3335
$ws.Range("J105").Value = 5197.8
$ws.Range("K105").Value = 4919.3335
$ws.Range("L105").Value = 5197.8
$ws.Range("M105").Value = -3172.3335
$ws.Range("N105").Value = -8691.799999999999

$ws.Range("H134").Value = 2539.6309
$ws.Range("I134").Value = 2272.2817
$ws.Range("J134").Value = 3999.7693
$ws.Range("K134").Value = 6816.8451
$ws.Range("L134").Value = 11999.3079
$ws.Range("M134").Value = -4281.8451
$ws.Range("N134").Value = -17069.3079

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4355.3447
$ws.Range("I31").Value = 3299.9443
$ws.Range("J31").Value = 6082.364
$ws.Range("K31").Value = 3299.9443
$ws.Range("L31").Value = 6082.364
$ws.Range("M31").Value = -3004.9443
$ws.Range("N31").Value = -6672.364

$ws.Range("H34").Value = 4355.3447
$ws.Range("I34").Value = 3299.9443
$ws.Range("J34").Value = 6082.364
$ws.Range("K34").Value = 3299.9443
$ws.Range("L34").Value = 6082.364
$ws.Range("M34").Value = -3097.9443
$ws.Range("N34").Value = -6486.364

$ws.Range("H58").Value = 37421.31
$ws.Range("I58").Value = 46409.652
$ws.Range("J58").Value = 2966
$ws.Range("K58").Value = 46409.652
$ws.Range("L58").Value = 2966
$ws.Range("M58").Value = -46206.652
$ws.Range("N58").Value = -3372

$ws.Range("H94").Value = 3224.3333
$ws.Range("I94").Value = 2424
$ws.Range("J94").Value = 3624.5
$ws.Range("K94").Value = 2424
$ws.Range("L94").Value = 3624.5
$ws.Range("M94").Value = -1973
$ws.Range("N94").Value = -4526.5

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws.Range("H99").Value = 3193.8572
$ws.Range("I99").Value = 2201.7778
$ws.Range("K99").Value = 2201.7778
$ws.Range("M99").Value = -703.7777999999998

$ws.Range("H105").Value = 2836.8
$ws.Range("I105").Value = 1705.5
$ws.Range("K105").Value = 1705.5
$ws.Range("M105").Value = 41.5

$ws.Range("H126").Value = 3193.8572
$ws.Range("I126").Value = 2201.7778
$ws.Range("K126").Value = 6605.3334
$ws.Range("M126").Value = -4135.3334

$ws.Range("H136").Value = 37421.31
$ws.Range("I136").Value = 46409.652
$ws.Range("J136").Value = 2966
$ws.Range("K136").Value = 139228.956
$ws.Range("L136").Value = 8898
$ws.Range("M136").Value = -136678.956
$ws.Range("N136").Value = -13998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 81.59999999999999
$ws.Range("I33").Value = 41.333332
$ws.Range("J33").Value = 142
$ws.Range("K33").Value = 247.999992
$ws.Range("L33").Value = 852
$ws.Range("M33").Value = 35.00000800000001
$ws.Range("N33").Value = -1418

$ws.Range("H68").Value = 728.1818
$ws.Range("I68").Value = 746
$ws.Range("J68").Value = 550
$ws.Range("K68").Value = 2238
$ws.Range("L68").Value = 1650
$ws.Range("M68").Value = -1427
$ws.Range("N68").Value = -3272

$ws.Range("H71").Value = 728.1818
$ws.Range("I71").Value = 746
$ws.Range("J71").Value = 550
$ws.Range("K71").Value = 6714
$ws.Range("L71").Value = 4950
$ws.Range("M71").Value = -2658
$ws.Range("N71").Value = -13062

$ws.Range("H75").Value = 733.3333
$ws.Range("J75").Value = 733.3333
$ws.Range("L75").Value = 2199.9999
$ws.Range("N75").Value = -4195.9999

$ws.Range("H78").Value = 733.3333
$ws.Range("J78").Value = 733.3333
$ws.Range("L78").Value = 6599.9997
$ws.Range("N78").Value = -16583.9997

$ws.Range("H107").Value = 373.38095
$ws.Range("I107").Value = 302.22223
$ws.Range("J107").Value = 426.75
$ws.Range("K107").Value = 906.66669
$ws.Range("L107").Value = 1280.25
$ws.Range("M107").Value = 1013.33331
$ws.Range("N107").Value = -5120.25

$ws.Range("H116").Value = 795.8333
$ws.Range("I116").Value = 795.8333
$ws.Range("K116").Value = 2387.4999
$ws.Range("M116").Value = 1054.5001

$ws.Range("H131").Value = 3855094.8
$ws.Range("J131").Value = 5567698.5
$ws.Range("L131").Value = 16703095.5
$ws.Range("N131").Value = -16713175.5

$ws.Range("H134").Value = 2281.0908
$ws.Range("J134").Value = 4259
$ws.Range("L134").Value = 12777
$ws.Range("N134").Value = -22917

$ws.Range("H137").Value = 4020.35
$ws.Range("I137").Value = 1586.1111
$ws.Range("J137").Value = 6012
$ws.Range("K137").Value = 4758.3333
$ws.Range("L137").Value = 18036
$ws.Range("M137").Value = 341.6666999999998
$ws.Range("N137").Value = -28236

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4664.5
$ws.Range("J80").Value = 4746
$ws.Range("L80").Value = 4746
$ws.Range("N80").Value = -6742

$ws.Range("H83").Value = 4664.5
$ws.Range("J83").Value = 4746
$ws.Range("L83").Value = 23730
$ws.Range("N83").Value = -33714

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws.Range("H132").Value = 41990
$ws.Range("I132").Value = 50311.477
$ws.Range("K132").Value = 150934.431
$ws.Range("M132").Value = -148404.431

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 6000
$ws.Range("J2").Value = 6000
$ws.Range("L2").Value = 6000
$ws.Range("N2").Value = -6224

$ws.Range("H7").Value = 10096.608
$ws.Range("I7").Value = 12968.8
$ws.Range("K7").Value = 12968.8
$ws.Range("M7").Value = -12856.8

$ws.Range("H22").Value = 70707.375
$ws.Range("J22").Value = 2727.8
$ws.Range("L22").Value = 2727.8
$ws.Range("N22").Value = -3317.8

$ws.Range("H27").Value = 70707.375
$ws.Range("J27").Value = 2727.8
$ws.Range("L27").Value = 2727.8
$ws.Range("N27").Value = -2941.8

$ws.Range("H40").Value = 3938
$ws.Range("I40").Value = 1951.4
$ws.Range("K40").Value = 1951.4
$ws.Range("M40").Value = -1815.4

$ws.Range("H43").Value = 75000
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()

$ws.Range("H61").Value = 3692.3333
$ws.Range("I61").Value = 3603.889
$ws.Range("J61").Value = 3957.6667
$ws.Range("K61").Value = 3603.889
$ws.Range("L61").Value = 3957.6667
$ws.Range("M61").Value = -3401.889
$ws.Range("N61").Value = -4361.6667

$ws.Range("H113").Value = 3692.3333
$ws.Range("I113").Value = 3603.889
$ws.Range("J113").Value = 3957.6667
$ws.Range("K113").Value = 3603.889
$ws.Range("L113").Value = 3957.6667
$ws.Range("M113").Value = -1433.889
$ws.Range("N113").Value = -8297.6667

$ws.Range("H122").Value = 3580.353
$ws.Range("I122").Value = 2776.3
$ws.Range("J122").Value = 4729
$ws.Range("K122").Value = 8328.900000000001
$ws.Range("L122").Value = 14187
$ws.Range("M122").Value = -5878.900000000001
$ws.Range("N122").Value = -19087

$ws.Range("H126").Value = 10096.608
$ws.Range("I126").Value = 12968.8
$ws.Range("K126").Value = 38906.39999999999
$ws.Range("M126").Value = -36436.39999999999

$ws.Range("H136").Value = 3526.6365
$ws.Range("I136").Value = 3526.6365
$ws.Range("K136").Value = 10579.9095
$ws.Range("M136").Value = -8029.9095

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 57747.168
$ws.Range("I126").Value = 113665.555
$ws.Range("K126").Value = 340996.665
$ws.Range("M126").Value = -338526.665
